$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Copy-Format($srcAddr, $dstAddr) {
    $ws.Range($srcAddr).Copy()
    $ws.Range($dstAddr).PasteSpecial(-4122)
}

# ============================================================
# Step 1: duplicate the existing "8,10,12,14_*" blocks (rows
# 1-61) down into rows 65-125, preserving their original
# labels/values exactly as they are before any further edits.
# ============================================================

# "8,10,12,14_randread_4k"  (rows 1-5 -> rows 65-69)
$ws.Range("B65:E65").Merge()
$ws.Range("B65").Value = "8,10,12,14_randread_4k"
Copy-Format "B1" "B65"

$ws.Range("B66").Value = 1
Copy-Format "B2" "B66"
$ws.Range("C66").Value = 2
Copy-Format "C2" "C66"
$ws.Range("D66").Value = 3
Copy-Format "D2" "D66"
$ws.Range("E66").Value = 4
Copy-Format "E2" "E66"

$ws.Range("A67").Value = "IOPS"
Copy-Format "A3" "A67"
$ws.Range("B67").Value = 6324
Copy-Format "B3" "B67"
$ws.Range("C67").Value = 11800
Copy-Format "C3" "C67"
$ws.Range("D67").Value = 16900
Copy-Format "D3" "D67"
$ws.Range("E67").Value = 22300
Copy-Format "E3" "E67"

$ws.Range("A68").Value = "BW(MB/s)"
Copy-Format "A4" "A68"
$ws.Range("B68").Value = 25.9
Copy-Format "B4" "B68"
$ws.Range("C68").Value = 48.1
Copy-Format "C4" "C68"
$ws.Range("D68").Value = 69.4
Copy-Format "D4" "D68"
$ws.Range("E68").Value = 91.4
Copy-Format "E4" "E68"

$ws.Range("A69").Value = "lat_avg"
Copy-Format "A5" "A69"
$ws.Range("B69").Value = 157.71748
Copy-Format "B5" "B69"
$ws.Range("C69").Value = 169.04784
Copy-Format "C5" "C69"
$ws.Range("D69").Value = 175.80612
Copy-Format "D5" "D69"
$ws.Range("E69").Value = 177.18539
Copy-Format "E5" "E69"

# "8,10,12,14_randread_128k"  (rows 9-13 -> rows 73-77)
$ws.Range("B73:E73").Merge()
$ws.Range("B73").Value = "8,10,12,14_randread_128k"
Copy-Format "B9" "B73"

$ws.Range("B74").Value = 1
Copy-Format "B10" "B74"
$ws.Range("C74").Value = 2
Copy-Format "C10" "C74"
$ws.Range("D74").Value = 3
Copy-Format "D10" "D74"
$ws.Range("E74").Value = 4
Copy-Format "E10" "E74"

$ws.Range("A75").Value = "IOPS"
Copy-Format "A11" "A75"
$ws.Range("B75").Value = 2820
Copy-Format "B11" "B75"
$ws.Range("C75").Value = 4785
Copy-Format "C11" "C75"
$ws.Range("D75").Value = 7228
Copy-Format "D11" "D75"
$ws.Range("E75").Value = 8770
Copy-Format "E11" "E75"

$ws.Range("A76").Value = "BW(MB/s)"
Copy-Format "A12" "A76"
$ws.Range("B76").Value = 370
Copy-Format "B12" "B76"
$ws.Range("C76").Value = 627
Copy-Format "C12" "C76"
$ws.Range("D76").Value = 947
Copy-Format "D12" "D76"
$ws.Range("E76").Value = 1150
Copy-Format "E12" "E76"

$ws.Range("A77").Value = "lat_avg"
Copy-Format "A13" "A77"
$ws.Range("B77").Value = 353.96
Copy-Format "B13" "B77"
$ws.Range("C77").Value = 408.66
Copy-Format "C13" "C77"
$ws.Range("D77").Value = 401.57
Copy-Format "D13" "D77"
$ws.Range("E77").Value = 448.9
Copy-Format "E13" "E77"

# "8,10,12,14_randwrite_4k"  (rows 17-21 -> rows 81-85)
$ws.Range("B81:E81").Merge()
$ws.Range("B81").Value = "8,10,12,14_randwrite_4k"
Copy-Format "B17" "B81"

$ws.Range("B82").Value = 1
Copy-Format "B18" "B82"
$ws.Range("C82").Value = 2
Copy-Format "C18" "C82"
$ws.Range("D82").Value = 3
Copy-Format "D18" "D82"
$ws.Range("E82").Value = 4
Copy-Format "E18" "E82"

$ws.Range("A83").Value = "IOPS"
Copy-Format "A19" "A83"
$ws.Range("B83").Value = 57300
Copy-Format "B19" "B83"
$ws.Range("C83").Value = 61700
Copy-Format "C19" "C83"
$ws.Range("D83").Value = 86300
Copy-Format "D19" "D83"
$ws.Range("E83").Value = 88600
Copy-Format "E19" "E83"

$ws.Range("A84").Value = "BW(MB/s)"
Copy-Format "A20" "A84"
$ws.Range("B84").Value = 235
Copy-Format "B20" "B84"
$ws.Range("C84").Value = 253
Copy-Format "C20" "C84"
$ws.Range("D84").Value = 354
Copy-Format "D20" "D84"
$ws.Range("E84").Value = 363
Copy-Format "E20" "E84"

$ws.Range("A85").Value = "lat_avg"
Copy-Format "A21" "A85"
$ws.Range("B85").Value = 12.02
Copy-Format "B21" "B85"
$ws.Range("C85").Value = 14.03
Copy-Format "C21" "C85"
$ws.Range("D85").Value = 11.31
Copy-Format "D21" "D85"
$ws.Range("E85").Value = 12.26
Copy-Format "E21" "E85"

# "8,10,12,14_randwrite_128k"  (rows 25-29 -> rows 89-93)
$ws.Range("B89:E89").Merge()
$ws.Range("B89").Value = "8,10,12,14_randwrite_128k"
Copy-Format "B25" "B89"

$ws.Range("B90").Value = 1
Copy-Format "B26" "B90"
$ws.Range("C90").Value = 2
Copy-Format "C26" "C90"
$ws.Range("D90").Value = 3
Copy-Format "D26" "D90"
$ws.Range("E90").Value = 4
Copy-Format "E26" "E90"

$ws.Range("A91").Value = "IOPS"
Copy-Format "A27" "A91"
$ws.Range("B91").Value = 1984
Copy-Format "B27" "B91"
$ws.Range("C91").Value = 2209
Copy-Format "C27" "C91"
$ws.Range("D91").Value = 3038
Copy-Format "D27" "D91"
$ws.Range("E91").Value = 3072
Copy-Format "E27" "E91"

$ws.Range("A92").Value = "BW(MB/s)"
Copy-Format "A28" "A92"
$ws.Range("B92").Value = 260
Copy-Format "B28" "B92"
$ws.Range("C92").Value = 290
Copy-Format "C28" "C92"
$ws.Range("D92").Value = 398
Copy-Format "D28" "D92"
$ws.Range("E92").Value = 403
Copy-Format "E28" "E92"

$ws.Range("A93").Value = "lat_avg"
Copy-Format "A29" "A93"
$ws.Range("B93").Value = 120.41
Copy-Format "B29" "B93"
$ws.Range("C93").Value = 134.65
Copy-Format "C29" "C93"
$ws.Range("D93").Value = 143.05
Copy-Format "D29" "D93"
$ws.Range("E93").Value = 146.41
Copy-Format "E29" "E93"

# "8,10,12,14_read_4k"  (rows 33-37 -> rows 97-101)
$ws.Range("B97:E97").Merge()
$ws.Range("B97").Value = "8,10,12,14_read_4k"
Copy-Format "B33" "B97"

$ws.Range("B98").Value = 1
Copy-Format "B34" "B98"
$ws.Range("C98").Value = 2
Copy-Format "C34" "C98"
$ws.Range("D98").Value = 3
Copy-Format "D34" "D98"
$ws.Range("E98").Value = 4
Copy-Format "E34" "E98"

$ws.Range("A99").Value = "IOPS"
Copy-Format "A35" "A99"
$ws.Range("B99").Value = 197000
Copy-Format "B35" "B99"
$ws.Range("C99").Value = 323000
Copy-Format "C35" "C99"
$ws.Range("D99").Value = 352000
Copy-Format "D35" "D99"
$ws.Range("E99").Value = 428000
Copy-Format "E35" "E99"

$ws.Range("A100").Value = "BW(MB/s)"
Copy-Format "A36" "A100"
$ws.Range("B100").Value = 809
Copy-Format "B36" "B100"
$ws.Range("C100").Value = 1322
Copy-Format "C36" "C100"
$ws.Range("D100").Value = 1443
Copy-Format "D36" "D100"
$ws.Range("E100").Value = 1754
Copy-Format "E36" "E100"

$ws.Range("A101").Value = "lat_avg"
Copy-Format "A37" "A101"
$ws.Range("B101").Value = 4.92645
Copy-Format "B37" "B101"
$ws.Range("C101").Value = 5.7801
Copy-Format "C37" "C101"
$ws.Range("D101").Value = 8.12196
Copy-Format "D37" "D101"
$ws.Range("E101").Value = 8.971969999999999
Copy-Format "E37" "E101"

# "8,10,12,14_read_128k"  (rows 41-45 -> rows 105-109)
$ws.Range("B105:E105").Merge()
$ws.Range("B105").Value = "8,10,12,14_read_128k"
Copy-Format "B41" "B105"

$ws.Range("B106").Value = 1
Copy-Format "B42" "B106"
$ws.Range("C106").Value = 2
Copy-Format "C42" "C106"
$ws.Range("D106").Value = 3
Copy-Format "D42" "D106"
$ws.Range("E106").Value = 4
Copy-Format "E42" "E106"

$ws.Range("A107").Value = "IOPS"
Copy-Format "A43" "A107"
$ws.Range("B107").Value = 4718
Copy-Format "B43" "B107"
$ws.Range("C107").Value = 8062
Copy-Format "C43" "C107"
$ws.Range("D107").Value = 12500
Copy-Format "D43" "D107"
$ws.Range("E107").Value = 13600
Copy-Format "E43" "E107"

$ws.Range("A108").Value = "BW(MB/s)"
Copy-Format "A44" "A108"
$ws.Range("B108").Value = 619
Copy-Format "B44" "B108"
$ws.Range("C108").Value = 1057
Copy-Format "C44" "C108"
$ws.Range("D108").Value = 1643
Copy-Format "D44" "D108"
$ws.Range("E108").Value = 1784
Copy-Format "E44" "E108"

$ws.Range("A109").Value = "lat_avg"
Copy-Format "A45" "A109"
$ws.Range("B109").Value = 210.46
Copy-Format "B45" "B109"
$ws.Range("C109").Value = 237.17
Copy-Format "C45" "C109"
$ws.Range("D109").Value = 234.55
Copy-Format "D45" "D109"
$ws.Range("E109").Value = 285.8
Copy-Format "E45" "E109"

# "8,10,12,14_write_4k"  (rows 49-53 -> rows 113-117)
$ws.Range("B113:E113").Merge()
$ws.Range("B113").Value = "8,10,12,14_write_4k"
Copy-Format "B49" "B113"

$ws.Range("B114").Value = 1
Copy-Format "B50" "B114"
$ws.Range("C114").Value = 2
Copy-Format "C50" "C114"
$ws.Range("D114").Value = 3
Copy-Format "D50" "D114"
$ws.Range("E114").Value = 4
Copy-Format "E50" "E114"

$ws.Range("A115").Value = "IOPS"
Copy-Format "A51" "A115"
$ws.Range("B115").Value = 66100
Copy-Format "B51" "B115"
$ws.Range("C115").Value = 121000
Copy-Format "C51" "C115"
$ws.Range("D115").Value = 169000
Copy-Format "D51" "D115"
$ws.Range("E115").Value = 200000
Copy-Format "E51" "E115"

$ws.Range("A116").Value = "BW(MB/s)"
Copy-Format "A52" "A116"
$ws.Range("B116").Value = 271
Copy-Format "B52" "B116"
$ws.Range("C116").Value = 496
Copy-Format "C52" "C116"
$ws.Range("D116").Value = 693
Copy-Format "D52" "D116"
$ws.Range("E116").Value = 821
Copy-Format "E52" "E116"

$ws.Range("A117").Value = "lat_avg"
Copy-Format "A53" "A117"
$ws.Range("B117").Value = 8.85
Copy-Format "B53" "B117"
$ws.Range("C117").Value = 8.74
Copy-Format "C53" "C117"
$ws.Range("D117").Value = 9.4
Copy-Format "D53" "D117"
$ws.Range("E117").Value = 9.94
Copy-Format "E53" "E117"

# "8,10,12,14_write_128k"  (rows 57-61 -> rows 121-125)
$ws.Range("B121:E121").Merge()
$ws.Range("B121").Value = "8,10,12,14_write_128k"
Copy-Format "B57" "B121"

$ws.Range("B122").Value = 1
Copy-Format "B58" "B122"
$ws.Range("C122").Value = 2
Copy-Format "C58" "C122"
$ws.Range("D122").Value = 3
Copy-Format "D58" "D122"
$ws.Range("E122").Value = 4
Copy-Format "E58" "E122"

$ws.Range("A123").Value = "IOPS"
Copy-Format "A59" "A123"
$ws.Range("B123").Value = 3056
Copy-Format "B59" "B123"
$ws.Range("C123").Value = 6627
Copy-Format "C59" "C123"
$ws.Range("D123").Value = 8347
Copy-Format "D59" "D123"
$ws.Range("E123").Value = 9287
Copy-Format "E59" "E123"

$ws.Range("A124").Value = "BW(MB/s)"
Copy-Format "A60" "A124"
$ws.Range("B124").Value = 401
Copy-Format "B60" "B124"
$ws.Range("C124").Value = 869
Copy-Format "C60" "C124"
$ws.Range("D124").Value = 1094
Copy-Format "D60" "D124"
$ws.Range("E124").Value = 1217
Copy-Format "E60" "E124"

$ws.Range("A125").Value = "lat_avg"
Copy-Format "A61" "A125"
$ws.Range("B125").Value = 92.4
Copy-Format "B61" "B125"
$ws.Range("C125").Value = 94.29
Copy-Format "C61" "C125"
$ws.Range("D125").Value = 100.65
Copy-Format "D61" "D125"
$ws.Range("E125").Value = 106.61
Copy-Format "E61" "E125"

# ============================================================
# Step 2: overwrite rows 1-61 with the new "0,2,4,6_*" results
# (styles/merges already in place - only labels + values change)
# ============================================================

# rows 1-5: "0,2,4,6_randread_4k"
$ws.Range("B1").Value = "0,2,4,6_randread_4k"
$ws.Range("B3").Value = 6247
$ws.Range("C3").Value = 12900
$ws.Range("D3").Value = 17700
$ws.Range("E3").Value = 23500
$ws.Range("B4").Value = 25.6
$ws.Range("C4").Value = 52.9
$ws.Range("D4").Value = 72.4
$ws.Range("E4").Value = 96.2
$ws.Range("B5").Value = 159.6543
$ws.Range("C5").Value = 153.38775
$ws.Range("D5").Value = 168.73814
$ws.Range("E5").Value = 169.44949

# rows 9-13: "0,2,4,6_randread_128k"
$ws.Range("B9").Value = "0,2,4,6_randread_128k"
$ws.Range("B11").Value = 4266
$ws.Range("C11").Value = 7236
$ws.Range("D11").Value = 12800
$ws.Range("E11").Value = 11800
$ws.Range("B12").Value = 559
$ws.Range("C12").Value = 949
$ws.Range("D12").Value = 1678
$ws.Range("E12").Value = 1543
$ws.Range("B13").Value = 233.29
$ws.Range("C13").Value = 270.66
$ws.Range("D13").Value = 222.53
$ws.Range("E13").Value = 330.77

# rows 17-21: "0,2,4,6_randwrite_4k"
$ws.Range("B17").Value = "0,2,4,6_randwrite_4k"
$ws.Range("B19").Value = 19600
$ws.Range("C19").Value = 36100
$ws.Range("D19").Value = 55600
$ws.Range("E19").Value = 56800
$ws.Range("B20").Value = 80.4
$ws.Range("C20").Value = 148
$ws.Range("D20").Value = 228
$ws.Range("E20").Value = 233
$ws.Range("B21").Value = 46.31
$ws.Range("C21").Value = 33.76
$ws.Range("D21").Value = 28.39
$ws.Range("E21").Value = 42.81

# rows 25-29: "0,2,4,6_randwrite_128k"
$ws.Range("B25").Value = "0,2,4,6_randwrite_128k"
$ws.Range("B27").Value = 1224
$ws.Range("C27").Value = 2235
$ws.Range("D27").Value = 3147
$ws.Range("E27").Value = 3138
$ws.Range("B28").Value = 161
$ws.Range("C28").Value = 293
$ws.Range("D28").Value = 413
$ws.Range("E28").Value = 411
$ws.Range("B29").Value = 337.48
$ws.Range("C29").Value = 228.05
$ws.Range("D29").Value = 196.47
$ws.Range("E29").Value = 285.35

# rows 33-37: "0,2,4,6_read_4k"
$ws.Range("B33").Value = "0,2,4,6_read_4k"
$ws.Range("B35").Value = 218000
$ws.Range("C35").Value = 279000
$ws.Range("D35").Value = 371000
$ws.Range("E35").Value = 749000
$ws.Range("B36").Value = 895
$ws.Range("C36").Value = 1142
$ws.Range("D36").Value = 1519
$ws.Range("E36").Value = 3068
$ws.Range("B37").Value = 4.46817
$ws.Range("C37").Value = 5.89138
$ws.Range("D37").Value = 7.85371
$ws.Range("E37").Value = 4.95728

# rows 41-45: "0,2,4,6_read_128k"
$ws.Range("B41").Value = "0,2,4,6_read_128k"
$ws.Range("B43").Value = 12000
$ws.Range("C43").Value = 23800
$ws.Range("D43").Value = 18000
$ws.Range("E43").Value = 16700
$ws.Range("B44").Value = 1699
$ws.Range("C44").Value = 3121
$ws.Range("D44").Value = 2486
$ws.Range("E44").Value = 2182
$ws.Range("B45").Value = 75.87
$ws.Range("C45").Value = 81.62
$ws.Range("D45").Value = 156.11
$ws.Range("E45").Value = 235.38

# rows 49-53: "0,2,4,6_write_4k"
$ws.Range("B49").Value = "0,2,4,6_write_4k"
$ws.Range("B51").Value = 57200
$ws.Range("C51").Value = 101000
$ws.Range("D51").Value = 144000
$ws.Range("E51").Value = 176000
$ws.Range("B52").Value = 234
$ws.Range("C52").Value = 415
$ws.Range("D52").Value = 590
$ws.Range("E52").Value = 720
$ws.Range("B53").Value = 10.96
$ws.Range("C53").Value = 11.31
$ws.Range("D53").Value = 11.59
$ws.Range("E53").Value = 11.88

# rows 57-61: "0,2,4,6_write_128k"
$ws.Range("B57").Value = "0,2,4,6_write_128k"
$ws.Range("B59").Value = 2782
$ws.Range("C59").Value = 6671
$ws.Range("D59").Value = 8084
$ws.Range("E59").Value = 8641
$ws.Range("B60").Value = 365
$ws.Range("C60").Value = 874
$ws.Range("D60").Value = 1060
$ws.Range("E60").Value = 1133
$ws.Range("B61").Value = 92.2
$ws.Range("C61").Value = 92.4
$ws.Range("D61").Value = 93.33
$ws.Range("E61").Value = 96.67
